# Adds verification-matrix rows for the new Fig. 5 plotting code:
#  - A new R script that plots the MetAtlas t-SNE (row 16, a previously
#    blank row right under the other "trivial figure generation" scripts)
#  - A new pair of "...8" GTEx-individual model generation scripts,
#    inserted right next to their existing (non-"8") counterparts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 16 was an empty row already sitting between the R-script block
#    (rows 2-15) and the MATLAB-script block (rows 18+) - just fill it in.
$ws.Range("A16").Value = "PlotMetAtlasTsne.R"
$ws.Range("B16").Value = "Not tested, only trivial figure generation code, the data is loaded from file."

# 2) Insert a new row right after "generate_gtexind_models.m" (row 57) for
#    the new "generate_gtexind_models8.m" script, pushing everything below
#    down by one.
$ws.Rows("58:58").Insert()
$ws.Range("A58").Value = "generate_gtexind_models8.m"
$ws.Range("B58").Value = "The code just calls model generation, it is not explicitly tested."

# 3) Insert another new row right after "gen_gtex_ind.sh" (now row 59) for
#    the new "gen_gtex_ind8.sh" script, pushing everything below down again.
$ws.Rows("60:60").Insert()
$ws.Range("A60").Value = "gen_gtex_ind8.sh"
$ws.Range("B60").Value = "trivial slurm script"

# Match the author's final view state: scrolled down a bit with B16 selected.
$ws.Range("B16").Select()
$excel.ActiveWindow.ScrollRow = 4
